$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target values for columns B-H (Total Item, Nill, Super Under Stock, Under Stock,
# Normal Stock, Over Stock, Super Over Stock) for rows 2-32.
$data = @{
    2 = @(51,0,45,5,1,0,0)
    3 = @(51,0,33,0,1,0,17)
    4 = @(51,0,20,0,0,0,31)
    5 = @(51,0,24,7,1,8,11)
    6 = @(51,1,21,5,5,2,17)
    7 = @(51,0,21,1,1,5,23)
    8 = @(51,0,20,5,2,4,20)
    9 = @(51,0,24,9,1,5,12)
    10 = @(51,0,23,10,2,6,10)
    11 = @(51,0,20,6,1,0,24)
    12 = @(51,0,21,4,0,4,22)
    13 = @(51,0,23,10,1,2,15)
    14 = @(51,0,22,6,5,5,13)
    15 = @(51,0,21,2,2,3,23)
    16 = @(51,0,24,7,3,5,12)
    17 = @(51,0,23,9,4,3,12)
    18 = @(51,0,22,7,4,4,14)
    19 = @(51,0,20,4,2,2,23)
    20 = @(51,0,21,2,3,6,19)
    21 = @(51,0,23,4,1,1,22)
    22 = @(51,0,27,5,3,3,13)
    23 = @(51,0,21,1,4,6,19)
    24 = @(51,0,26,5,4,3,13)
    25 = @(51,0,21,3,3,3,21)
    26 = @(51,0,25,6,11,3,6)
    27 = @(51,0,22,3,2,5,19)
    28 = @(51,0,21,0,0,1,29)
    29 = @(51,0,22,3,4,4,18)
    30 = @(51,0,21,6,6,3,15)
    31 = @(51,0,22,8,2,3,16)
    32 = @(51,0,22,9,5,5,10)
}

foreach ($r in $data.Keys) {
    $vals = $data[$r]
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $col = 2 + $i   # column B = 2
        $ws.Cells.Item($r, $col).Value = $vals[$i]
    }
}
